$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Approbateurs" column, row "Alban PAPASSIAN" -> add "Fanny LAJEUNESSE"
#    in the (currently empty) approver cell of the 1st table.
# ------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$approverCell = $t1.Cell(2, 2)
$approverCell.Range.InsertAfter("Fanny LAJEUNESSE")

# Re-fetch the cell/range after the mutation and apply the same font
# size (11.5pt = w:sz 23 half-points) used throughout this table so the
# new run carries <w:rPr><w:sz w:val="23"/></w:rPr>, matching its
# neighbours.
$t1b = $d.Tables.Item(1)
$approverCell2 = $t1b.Cell(2, 2)
$approverRange = $approverCell2.Range
$approverRange.Font.Size = 11.5

# ------------------------------------------------------------------
# 2) "Approuvé le :" cell -> append the approval date "05/11" right
#    after the existing trailing space run, then drop the (moved)
#    "_GoBack" bookmark immediately after the new text.
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(1)
$dateCell = $t2.Cell(3, 2)
# Insert the date plus one throwaway sentinel character. Placing the
# bookmark touching the very last character of the very last run in a
# paragraph is mishandled when done directly, so the sentinel keeps the
# bookmark's anchor away from that edge case; it is deleted afterwards.
$dateCell.Range.InsertAfter("05/11X")

$t2b = $d.Tables.Item(1)
$dateCell2 = $t2b.Cell(3, 2)
$dateCellRange = $dateCell2.Range
# The last 7 positions up to (and including) the cell mark are, in
# order: '0' '5' '/' '1' '1' 'X' <cell-mark>; style just "05/11".
$dateStart = $dateCellRange.End - 7
$newDateRange = $d.Range($dateStart, $dateStart + 5)
$newDateRange.Font.Size = 11.5

# Bookmark goes right after "05/11" (i.e. right before the sentinel "X").
$bookmarkPos = $dateStart + 5
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the sentinel "X" now that the bookmark is anchored correctly
# (the bookmark itself is zero-width, so "X" still starts at $bookmarkPos).
$sentinelRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinelRange.Delete()
